$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value2 = 1.448123574256897
$ws.Range("B1").Value2 = 2.342830419540405
$ws.Range("C1").Value2 = 2.908051013946533
$ws.Range("D1").Value2 = 3.377647876739502
$ws.Range("E1").Value2 = 2.079902172088623
